# Price-tracker update: append a new row for the 2026-02-07 scrape.
# The new row repeats the Price/Discount/Incredible values of the
# previous (most recent) row, only the Date column changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Locate the last used row in column A (the Date column) and the row
# right after it, where the new entry will go.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$srcRange = $ws.Range("A" + $lastRow + ":D" + $lastRow)
$dstRange = $ws.Range("A" + $newRow + ":D" + $newRow)

# Duplicate the previous row (values + default styling) into the new
# row so Price/Discount/Incredible come along unchanged and no new
# cell styles are introduced.
$srcRange.Copy($dstRange)

# Overwrite just the Date cell with the new scrape date. The leading
# apostrophe keeps Excel from reinterpreting the text as a date
# serial number, so it is stored as a plain string like the rest of
# the column; re-applying the Normal style afterwards clears the
# "quote prefix" formatting Excel adds for that apostrophe.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = "'2026-02-07"
$dateCell.Style = "Normal"

$wb.Save()
